$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3: update URL text, add hyperlink, and match the existing hyperlink cell style (same as D2/D5)
$ws.Range("D3").Value = "https://masterdaily.dev.webchart.app/webchart.cgi"
$ws.Hyperlinks.Add($ws.Range("D3"), "https://masterdaily.dev.webchart.app/webchart.cgi")
$ws.Range("D3").Style = $ws.Range("D2").Style

# D4: update hyperlink target URL text (keeps its existing style/hyperlink)
$ws.Range("D4").Value = "https://masterdaily.dev.webchart.app/webchart.cgi"

# C4: update prescriber + patient name in the instructions text
$ws.Range("C4").Value = 'Write Amoxicillin 500mg capsule 2 caps daily for 7 days. For Prescriber: Your name. Total quantity: 14 and no refills. Allow substitutions Send the script to "MIE Test Pharmacy", for patient HART, WILLIAM S'

# E4: update verification text to match new patient name
$ws.Range("E4").Value = "verify text equals `"HART, WILLIAM S.`" in `"patient_name`"`nverify text equals `"amoxicillin 500 mg tablet`" in `"medicine`""

# Move the active selection to D3, matching the edited cell
$ws.Range("D3").Select()
